# Jinzu test case complete
# - Rename the three existing sheets
# - Duplicate the (renamed) "getDataGraphQL-iems" sheet into a new sheet named
#   "getDataGraphQL" that becomes the active tab, with new/extended test data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "getConceptModelDataByCondition-"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "getDataEntities-iems"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "getDataGraphQL-iems"

# Duplicate sheet3 right after itself; the copy becomes sheet index 4.
$ws3.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "getDataGraphQL"

# Extend the table with two more data rows (5 and 6), carrying the same
# cell formatting as the existing data rows (row 4's style).
$ws4.Range("A4:F4").Copy()
$ws4.Range("A5:F6").PasteSpecial(-4122)

# Row 4: var23 / AND query
$ws4.Range("A4").Value = "JinZu-ApiEngine-Test-7-var23"
$ws4.Range("C4").Value = '{Project (cond:"{_and: [{business_mgr:{_in:[\"潘云晖\",\"臧佳宝\" ]}},{status:{_in:[\"archived\" ]}}]}",order:"") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}'

# Row 5: var24 / OR query (new row)
$ws4.Range("A5").Value = "JinZu-ApiEngine-Test-7-var24"
$ws4.Range("C5").Value = '{Project (cond:"{_or: [{business_mgr:{_in:[\"潘云晖\",\"臧佳宝\" ]}},{status:{_in:[\"archived\" ]}}]}",order:"") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status}}'

# Row 6: var25 (new row)
$ws4.Range("A6").Value = "JinZu-ApiEngine-Test-7-var25"

# Row 3: var18
$ws4.Range("A3").Value = "JinZu-ApiEngine-Test-7-var18"

# Row 2: var19
$ws4.Range("A2").Value = "JinZu-ApiEngine-Test-7-var19"

# Row 2: business_mgr query
$ws4.Range("C2").Value = '{Project(cond:"{business_mgr:{_in:[\"潘云晖\",\"臧佳宝\" ]}}",order:"") {business_mgr business_unit charge_frequency city province district credit_amount detail_address discount_ratio expire_date guarantee_type id no status name risk_mgr rent_type}}'

# Row 3: status query
$ws4.Range("C3").Value = '{Project(cond:"{status:{_in:[\"online\",\"archived\" ]}}",order:"") {business_mgr business_unit charge_frequency city province district credit_amount detail_address discount_ratio expire_date guarantee_type id no status name risk_mgr rent_type}}'

# Row 6: OR + invert_Customer query
$ws4.Range("C6").Value = '{Project(cond:"{_or: [{business_mgr:{_in:[\"潘云晖\",\"臧佳宝\" ]}},{status:{_in:[\"archived\" ]}}]}",order:"") {business_mgr business_unit charge_frequency city class_level classification_level credit_amount detail_address discount_ratio district expire_date guarantee_type id is_manufacture_buy_back is_manufacture_leasing manufacture name no province rent_type risk_mgr status  invert_Customer (cond:"{id:{_eq:24}}",order:"") {actual_controller category cid city cname contact contact_detail ctype district enterprise_size group holding_type id is_connected_tx is_gov_fin_customer is_group_customer legal_person legal_person_id major_class middle_class office_address project province registered_address small_class}}}'

# Fix up remaining B/D/E/F columns so rows read "good request, data retrieved" / 200 / 100000 / Successfully
$ws4.Range("B2").Value = "good request, data retrieved"
$ws4.Range("E2").Value = 100000
$ws4.Range("F2").Value = "Successfully"

$ws4.Range("B5").Value = "good request, data retrieved"
$ws4.Range("D5").Value = 200
$ws4.Range("E5").Value = 100000
$ws4.Range("F5").Value = "Successfully"

$ws4.Range("B6").Value = "good request, data retrieved"
$ws4.Range("D6").Value = 200
$ws4.Range("E6").Value = 100000
$ws4.Range("F6").Value = "Successfully"
